# Applies DBC-export style updates to the "AMS" worksheet:
#  - cell_voltage_* signal rows: Factor 1 -> 0.001, Unit "" -> "V"
#  - temperature_value_* signal rows: Factor 1 -> 0.01, Unit "" -> "ºC"
#  - a handful of cell_voltage_* rows (121-132) only get the Factor update

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AMS")

# cell_voltage_* rows: Factor -> 0.001, Unit -> "V"
$ws.Range("F10:F13").Value = 0.001
$ws.Range("J10:J13").Value = "V"
$ws.Range("F17:F20").Value = 0.001
$ws.Range("J17:J20").Value = "V"
$ws.Range("F59:F62").Value = 0.001
$ws.Range("J59:J62").Value = "V"
$ws.Range("F66:F69").Value = 0.001
$ws.Range("J66:J69").Value = "V"
$ws.Range("F73:F76").Value = 0.001
$ws.Range("J73:J76").Value = "V"
$ws.Range("F115:F118").Value = 0.001
$ws.Range("J115:J118").Value = "V"
$ws.Range("F122:F125").Value = 0.001
$ws.Range("J122:J125").Value = "V"
$ws.Range("F129:F132").Value = 0.001
$ws.Range("J129:J132").Value = "V"
$ws.Range("F171:F174").Value = 0.001
$ws.Range("J171:J174").Value = "V"
$ws.Range("F178:F181").Value = 0.001
$ws.Range("J178:J181").Value = "V"
$ws.Range("F185:F188").Value = 0.001
$ws.Range("J185:J188").Value = "V"
$ws.Range("F227:F230").Value = 0.001
$ws.Range("J227:J230").Value = "V"
$ws.Range("F234:F237").Value = 0.001
$ws.Range("J234:J237").Value = "V"
$ws.Range("F241:F244").Value = 0.001
$ws.Range("J241:J244").Value = "V"
$ws.Range("F283:F286").Value = 0.001
$ws.Range("J283:J286").Value = "V"
$ws.Range("F290:F293").Value = 0.001
$ws.Range("J290:J293").Value = "V"
$ws.Range("F297:F300").Value = 0.001
$ws.Range("J297:J300").Value = "V"
$ws.Range("F339:F342").Value = 0.001
$ws.Range("J339:J342").Value = "V"
$ws.Range("F346:F349").Value = 0.001
$ws.Range("J346:J349").Value = "V"
$ws.Range("F353:F356").Value = 0.001
$ws.Range("J353:J356").Value = "V"
$ws.Range("F395:F398").Value = 0.001
$ws.Range("J395:J398").Value = "V"
$ws.Range("F402:F405").Value = 0.001
$ws.Range("J402:J405").Value = "V"
$ws.Range("F409:F412").Value = 0.001
$ws.Range("J409:J412").Value = "V"
$ws.Range("F451:F454").Value = 0.001
$ws.Range("J451:J454").Value = "V"
$ws.Range("F458:F461").Value = 0.001
$ws.Range("J458:J461").Value = "V"
$ws.Range("F465:F468").Value = 0.001
$ws.Range("J465:J468").Value = "V"
$ws.Range("F507:F510").Value = 0.001
$ws.Range("J507:J510").Value = "V"
$ws.Range("F514:F517").Value = 0.001
$ws.Range("J514:J517").Value = "V"
$ws.Range("F521:F524").Value = 0.001
$ws.Range("J521:J524").Value = "V"
$ws.Range("F619:F622").Value = 0.001
$ws.Range("J619:J622").Value = "V"
$ws.Range("F626:F629").Value = 0.001
$ws.Range("J626:J629").Value = "V"
$ws.Range("F633:F636").Value = 0.001
$ws.Range("J633:J636").Value = "V"

# temperature_value_* rows: Factor -> 0.01, Unit -> "ºC"
$ws.Range("F24:F27").Value = 0.01
$ws.Range("J24:J27").Value = "ºC"
$ws.Range("F31:F34").Value = 0.01
$ws.Range("J31:J34").Value = "ºC"
$ws.Range("F38:F41").Value = 0.01
$ws.Range("J38:J41").Value = "ºC"
$ws.Range("F80:F83").Value = 0.01
$ws.Range("J80:J83").Value = "ºC"

# cell_voltage_121..132 rows: Factor -> 0.001 only (Unit left unchanged)
$ws.Range("F563:F566").Value = 0.001
$ws.Range("F570:F573").Value = 0.001
$ws.Range("F577:F580").Value = 0.001
